$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.95"
$ws.Range("E2").Value = "'-3.29%"
$ws.Range("D3").Value = "'54.16"
$ws.Range("E3").Value = "'10.43%"
$ws.Range("D4").Value = "'5.094"
$ws.Range("E4").Value = "'-4.12%"
$ws.Range("D5").Value = "'0.07908"
$ws.Range("E5").Value = "'-1.89%"
$ws.Range("D6").Value = "'4.561"
$ws.Range("E6").Value = "'-0.94%"
$ws.Range("D7").Value = "'1.395"
$ws.Range("E7").Value = "'3.88%"
$ws.Range("D8").Value = "'1.671"
$ws.Range("E8").Value = "'1.74%"
$ws.Range("D9").Value = "'0.1244"
$ws.Range("E9").Value = "'-2.81%"
$ws.Range("D10").Value = "'0.2015"
$ws.Range("E10").Value = "'2.43%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09507"
$ws.Range("E11").Value = "'-1.38%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04718"
$ws.Range("E12").Value = "'-0.14%"
$ws.Range("D13").Value = "'0.1045"
$ws.Range("E13").Value = "'-0.19%"
$ws.Range("D14").Value = "'0.001275"
$ws.Range("E14").Value = "'-3.42%"
$ws.Range("D15").Value = "'0.005827"
$ws.Range("E15").Value = "'-1.56%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.345"
$ws.Range("E16").Value = "'-0.06%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.436"
$ws.Range("E17").Value = "'-0.32%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3429"
$ws.Range("E18").Value = "'-2.30%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'8.367"
$ws.Range("E19").Value = "'4.59%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1363"
$ws.Range("E20").Value = "'-0.17%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.2913"
$ws.Range("E21").Value = "'-5.73%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04171"
$ws.Range("E22").Value = "'-0.51%"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.001258"
$ws.Range("E23").Value = "'-4.17%"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.003988"
$ws.Range("E24").Value = "'-8.05%"
$ws.Range("B25").Value = "NitroEx"
$ws.Range("C25").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D25").Value = "'0.0001347"
$ws.Range("E25").Value = "'-0.04%"
$ws.Range("B26").Value = "UpBots"
$ws.Range("C26").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D26").Value = "'0.0003534"
$ws.Range("E26").Value = "'-0.02%"
$ws.Range("D38").Value = "'0.02645"
$ws.Range("E38").Value = "'-3.02%"
$ws.Range("D39").Value = "'0.05942"
$ws.Range("E39").Value = "'-0.54%"
$ws.Range("D40").Value = "'0.01080"
$ws.Range("E40").Value = "'-0.39%"
$ws.Range("D41").Value = "'0.1766"
$ws.Range("E41").Value = "'20.43%"
$ws.Range("D42").Value = "'0.007925"
$ws.Range("E42").Value = "'-0.97%"
$ws.Range("D43").Value = "'0.008186"
$ws.Range("E43").Value = "'3.74%"
$ws.Range("D44").Value = "'0.008374"
$ws.Range("E44").Value = "'6.17%"
$ws.Range("D45").Value = "'0.3411"
$ws.Range("E45").Value = "'-2.82%"
$ws.Range("D46").Value = "'0.00007177"
$ws.Range("E46").Value = "'4.14%"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("D48").Value = "'0.05532"
$ws.Range("E48").Value = "'-5.54%"
$ws.Range("D49").Value = "'0.002616"
$ws.Range("E49").Value = "'-34.52%"
$ws.Range("D50").Value = "'0.00002096"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("D51").Value = "'0.0001996"
$ws.Range("E51").Value = "'-0.06%"
